$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New metadata columns: data_collection_start_date, data_collection_end_date, publication_date ---
$ws.Range("K1").Value = "data_collection_start_date"
$ws.Range("L1").Value = "data_collection_end_date"
$ws.Range("M1").Value = "publication_date"

# Match the bold header style used by the rest of row 1 (e.g. D1/J1)
$ws.Range("K1:M1").Font.Bold = $true

# --- Fill in the previously-empty "name" value for this baseline row, with wrap text ---
$ws.Range("E2").Value = "Kasungu Lilongwe Plains "
$ws.Range("E2").WrapText = $true

# --- Column widths for the new / widened columns ---
$ws.Columns("E").ColumnWidth = 21.25
$ws.Columns("K").ColumnWidth = 23.35
$ws.Columns("L").ColumnWidth = 22.48
$ws.Columns("M").ColumnWidth = 14.9

# Selection moves to the newly-populated E2 cell
$ws.Range("E2").Select() | Out-Null

Write-Output "Applied metadata updates"
